$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-blank Percent column (C) for rows 25-34, matching
# the percentage number format already used by the rest of column C
# (e.g. C22, style index 4 -> numFmtId 9, "0%").
$fmt = $ws.Range("C22").NumberFormat

$values = @{
    25 = 0.09
    26 = 0.08
    27 = 0.08
    28 = 0.07
    29 = 0.07
    30 = 0.1
    31 = 0.09
    32 = 0.08
    33 = 0.08
    34 = 0.07
}

foreach ($row in 25..34) {
    $cell = $ws.Range("C$row")
    $cell.Value = $values[$row]
    $cell.NumberFormat = $fmt
}

# Update the visible selection to C2 (matches the saved view state in the
# workbook after upload).
[void]$ws.Range("C2").Select()
